# Generate Report for Handoff
#
# Rows 4-7 (347a8044..., 4ff9034b..., 8bea6420..., e75bda2e...) moved from
# "low" priority to "ht" (handed-off) priority, and their Latest Handoff
# Datetime got refreshed to reflect the new handoff pass, on both the
# zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4..7 -> Priority "ht", Latest Handoff Datetime refreshed
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-04 00:35:29"
}

# de-de sheet: rows 4..7 -> Priority "ht", Latest Handoff Datetime refreshed
foreach ($r in 4..7) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-09-04 00:35:34"
}
